{"js": "// Change the \"Integrantes:\" header name from \"Juan Gomez\" to\n// \"Juan Ignacio Gomez\" by splitting the single run that holds\n// \"Juan Gomez\" into three runs: \"Juan\", \" Ignacio\", \" Gomez\".\n\n// The name lives in the default page header, not in the main body,\n// so we have to reach it through the section's header story.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst section = sections.items[0];\nconst header = section.getHeader(\"Primary\");\n\n// Locate the run of text we need to touch.\nconst found = header.search(\"Juan Gomez\", { matchCase: true, matchWholeWord: false });\nfound.load(\"items\");\nawait context.sync();\n\nif (found.items.length === 0) {\n  throw new Error('Could not find \"Juan Gomez\" in the primary header.');\n}\n\nconst target = found.items[0];\n\n// Replace the found range with the equivalent OOXML: the same run\n// (keeping its original rsid) now just containing \"Juan\", followed\n// by two brand-new runs for \" Ignacio\" and \" Gomez\". Using insertOoxml\n// (instead of a plain text insertText, which Word/the host would fold\n// back into a single run) is what lets the result land as three\n// distinct <w:r> elements, matching how the document was actually\n// edited.\nconst flatOpcXml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r w:rsidR=\"009C7224\"><w:rPr><w:lang w:val=\"es-MX\"/></w:rPr><w:t>Juan</w:t></w:r>\n            <w:r><w:rPr><w:lang w:val=\"es-MX\"/></w:rPr><w:t xml:space=\"preserve\"> Ignacio</w:t></w:r>\n            <w:r><w:rPr><w:lang w:val=\"es-MX\"/></w:rPr><w:t xml:space=\"preserve\"> Gomez</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ntarget.insertOoxml(flatOpcXml, \"Replace\");\nawait context.sync();\n", "ps1": "# Change the \"Integrantes:\" header name from \"Juan Gomez\" to\n# \"Juan Ignacio Gomez\". In the saved OOXML this shows up as the single\n# run holding \"Juan Gomez\" turning into three runs: \"Juan\", \" Ignacio\"\n# and \" Gomez\" (same formatting throughout - <w:lang w:val=\"es-MX\"/>).\n\n$d = $word.ActiveDocument\n\n# The name lives in the default (primary) page header, not the body.\n$section = $d.Sections.Item(1)\n$header = $section.Headers.Item($wdHeaderFooterPrimary)\n\n# Find the exact text we need to touch.\n$rng = $header.Range.Duplicate\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"Juan Gomez\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$found = $find.Execute()\n\nif (-not $found) {\n    throw 'Could not find \"Juan Gomez\" in the primary header.'\n}\n\n# Replace the matched range's contents with equivalent WordprocessingML:\n# the original run (its rsid preserved) now holding just \"Juan\",\n# followed by two new runs for \" Ignacio\" and \" Gomez\". Plain\n# \"$rng.Text = ...\" / InsertAfter calls would have Word fold the new\n# characters back into a single run, so we use Range.InsertXML to land\n# three distinct <w:r> elements, matching how the document was\n# actually edited.\n$xml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:r w:rsidR=\"009C7224\"><w:rPr><w:lang w:val=\"es-MX\"/></w:rPr><w:t>Juan</w:t></w:r>' +\n    '<w:r><w:rPr><w:lang w:val=\"es-MX\"/></w:rPr><w:t xml:space=\"preserve\"> Ignacio</w:t></w:r>' +\n    '<w:r><w:rPr><w:lang w:val=\"es-MX\"/></w:rPr><w:t xml:space=\"preserve\"> Gomez</w:t></w:r>' +\n    '</w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n\n$rng.InsertXML($xml)\n"}
